$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.977.89"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "'2.690.76"
$ws.Range("E3").Value = "  +4.78%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'514.36"
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("D6").Value = "'144.31"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("D9").Value = "'2.686.67"
$ws.Range("E9").Value = "  +4.42%  "
$ws.Range("D10").Value = "'6.27"
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("E11").Value = "  +4.27%  "
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("E13").Value = "  -1.25%  "
$ws.Range("D14").Value = "'3.155.45"
$ws.Range("E14").Value = "  +4.71%  "
$ws.Range("D15").Value = "'58.957.91"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "'20.98"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").Value = "'2.681.45"
$ws.Range("E18").Value = "  +4.38%  "
$ws.Range("D19").Value = "'348.74"
$ws.Range("E19").Value = "  +4.54%  "
$ws.Range("D20").Value = "'4.54"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("D21").Value = "'10.42"
$ws.Range("E21").Value = "  +3.27%  "
$ws.Range("E22").Value = "  +2.82%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "'60.87"
$ws.Range("E24").Value = "  +1.83%  "
$ws.Range("E25").Value = "  +2.93%  "
$ws.Range("D26").Value = "'2.772.16"
$ws.Range("E26").Value = "  +3.69%  "
$ws.Range("D27").Value = "'0.995"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("D29").Value = "'0.0₃0809"
$ws.Range("E29").Value = "  +3.58%  "
$ws.Range("E30").Value = "  +4.71%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").Value = "'6.39"
$ws.Range("E32").Value = "  +9.24%  "
$ws.Range("D33").Value = "'18.89"
$ws.Range("E33").Value = "  +1.74%  "
$ws.Range("E34").Value = "  +1.96%  "
$ws.Range("D35").Value = "'149.90"
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("D36").Value = "'1.02"
$ws.Range("E36").Value = "  +13.54%  "
$ws.Range("E37").Value = "  +1.73%  "
$ws.Range("E38").Value = "  +3.26%  "
$ws.Range("D39").Value = "'36.75"
$ws.Range("E39").Value = "  +2.28%  "
$ws.Range("D40").Value = "'0.846"
$ws.Range("E40").Value = "  +2.65%  "
$ws.Range("D41").Value = "'3.67"
$ws.Range("E41").Value = "  +4.32%  "
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").Value = "'0.620"
$ws.Range("E43").Value = "  +1.92%  "
$ws.Range("D44").Value = "'278.76"
$ws.Range("E44").Value = "  -3.38%  "
$ws.Range("D45").Value = "'0.997"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("D47").Value = "'19.57"
$ws.Range("E47").Value = "  +4.23%  "
$ws.Range("E48").Value = "  -0.32%  "

# Row 49: RenderToken -> VeChain
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0230"
$ws.Range("E49").Value = "  +1.22%  "

# Row 50: VeChain -> WhiteBITCoin
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").Value = "'10.27"
$ws.Range("E50").Value = "  -0.62%  "

# Row 51: WhiteBITCoin -> Maker
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "'1.995.33"
$ws.Range("E51").Value = "  +4.08%  "
